$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new review text as a new row at the end (row 11)
$ws.Range("A11").Value = "очень хороший телефон, мне безумно понравился"
$ws.Range("A11").WrapText = $true

# Row heights auto-fit (wrap text) for the already-populated review rows
$ws.Range("A2").RowHeight = 409.6
$ws.Range("A3").RowHeight = 187.2
$ws.Range("A4").RowHeight = 360
$ws.Range("A5").RowHeight = 129.6
$ws.Range("A6").RowHeight = 259.2
$ws.Range("A7").RowHeight = 409.6
$ws.Range("A8").RowHeight = 409.6
$ws.Range("A9").RowHeight = 216
$ws.Range("A10").RowHeight = 409.6

# Move selection / view to reflect the new active cell / scrolled position
$ws.Range("A11").Select()
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("A14").Select()
